$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Commentaire" header in column L, matching the style used by
# the other header cells (copy format from K1, then set the text).
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("L1").Value = "Commentaire"

# Give the new column a sensible width (closest reachable value to the
# 13.29-character width used in the authored workbook).
$ws.Columns.Item(12).ColumnWidth = 12.5

# Reflect the cursor position recorded in the saved workbook.
$ws.Range("M5").Select() | Out-Null
